$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

$tbl = $ws.ListObjects.Item("Table3")

$row1 = $tbl.ListRows.Add()
$row2 = $tbl.ListRows.Add()

# Copy formatting (including borders/number formats) from the last existing data row (row 6)
# down onto the two freshly added rows so they match the table's established look.
$ws.Range("A6:F6").Copy() | Out-Null
$ws.Range("A7:F8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$row1.Range.Item(1,2).Value = "David SpiegelHalter"
$row1.Range.Item(1,1).Value = "The Art of Statistics"
$row1.Range.Item(1,3).Value2 = 45522
$row1.Range.Item(1,4).Value2 = 45609
$row1.Range.Item(1,5).Value = "***"
$row1.Range.Item(1,6).Value = "for a statistics book, it was pretty good. But still tough to get through frankly"

$row2.Range.Item(1,1).Value = "Energy Trading & Risk Management"
$row2.Range.Item(1,2).Value = "Steven Berley"
$row2.Range.Item(1,3).Value2 = 45547
$row2.Range.Item(1,4).Value2 = 45607
$row2.Range.Item(1,5).Value = "*"
$row2.Range.Item(1,6).Value = "this book was overly vague the whole time. I can't say I got a lot from it. "

$ws.Range("C7:D8").NumberFormat = "m/d/yyyy"

$ws.Range("A9").Select()
